# Add column 'name' in user table
# The "User" table occupies columns G (label), H (column name),
# I (data type) and J (notes) starting at row 1 (header row with
# G1=User, H1=Account, I1=varchar) followed by one table-row per
# user-table column in rows 2-4, plus two note rows (J5/J6 before
# the edit).
#
# A new row describing the 'name' column (type varchar) needs to be
# inserted right after the header, i.e. at row 2, pushing the
# existing rows (Email, user_pass, authority + its notes) down by
# one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the existing User-table rows (H2:J6) down one row, to H3:J7,
# working from the bottom up so values are not overwritten before
# being copied.
for ($r = 7; $r -ge 3; $r--) {
    $ws.Range("H" + $r).Value = $ws.Range("H" + ($r - 1)).Value2
    $ws.Range("I" + $r).Value = $ws.Range("I" + ($r - 1)).Value2
    $ws.Range("J" + $r).Value = $ws.Range("J" + ($r - 1)).Value2
}

# Insert the new 'name' / 'varchar' row for the user table right
# after the header row. (J2 has no note, so it is left blank.)
$ws.Range("H2").Value = "name"
$ws.Range("I2").Value = "varchar"
